$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update changed cells (ligand/receptor/target refs + recomputed TPM-based specificity values) ---
$ws.Range("B2").Value = "Nlgn1"
$ws.Range("C2").Value = "Nrxn2"
$ws.Range("D2").Value = "ECs"
$ws.Range("I2").Value = 0.1532734190501181
$ws.Range("J2").Value = 0.1532734190501181
$ws.Range("M2").Value = 0.003710666666666666
$ws.Range("N2").Value = 0.011132
$ws.Range("O2").Value = 0.001642024256586498
$ws.Range("P2").Value = 0.001642024256586498
$ws.Range("Q2").Value = 0.00006935978133333333
$ws.Range("R2").Value = 0.000624238032
$ws.Range("S2").Value = 0.000251678671970241
$ws.Range("T2").Value = 0.000251678671970241

# --- Row 3: update changed cells (ligand/receptor/target refs + recomputed TPM-based specificity values) ---
$ws.Range("B3").Value = "Nlgn1"
$ws.Range("C3").Value = "Nrxn2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("I3").Value = 0.1532734190501181
$ws.Range("J3").Value = 0.1532734190501181
$ws.Range("M3").Value = 2.049608666666666
$ws.Range("N3").Value = 6.148826
$ws.Range("O3").Value = 0.9069818039462568
$ws.Range("P3").Value = 0.9069818039462569
$ws.Range("Q3").Value = 0.03831128519733333
$ws.Range("S3").Value = 0.1390162021070867
$ws.Range("T3").Value = 0.1390162021070867

# --- Row 4: update changed cells (ligand/receptor/target refs + recomputed TPM-based specificity values) ---
$ws.Range("B4").Value = "Nlgn1"
$ws.Range("C4").Value = "Nrxn2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 0.1532734190501181
$ws.Range("J4").Value = 0.1532734190501181
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.206493
$ws.Range("N4").Value = 0.619479
$ws.Range("O4").Value = 0.09137617179715662
$ws.Range("P4").Value = 0.09137617179715661
$ws.Range("Q4").Value = 0.003859767156
$ws.Range("R4").Value = 0.034737904404
$ws.Range("S4").Value = 0.01400553827106117
$ws.Range("T4").Value = 0.01400553827106117

# --- Row 5 (new row) ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Nlgn1"
$ws.Range("C5").Value = "Nrxn2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.10326
$ws.Range("H5").Value = 0.30978
$ws.Range("I5").Value = 0.846726580949882
$ws.Range("J5").Value = 0.8467265809498818
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.003710666666666666
$ws.Range("N5").Value = 0.011132
$ws.Range("O5").Value = 0.001642024256586498
$ws.Range("P5").Value = 0.001642024256586498
$ws.Range("Q5").Value = 0.00038316344
$ws.Range("R5").Value = 0.00344847096
$ws.Range("S5").Value = 0.001390345584616257
$ws.Range("T5").Value = 0.001390345584616257

# --- Row 6 (new row) ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Nlgn1"
$ws.Range("C6").Value = "Nrxn2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.10326
$ws.Range("H6").Value = 0.30978
$ws.Range("I6").Value = 0.846726580949882
$ws.Range("J6").Value = 0.8467265809498818
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.049608666666666
$ws.Range("N6").Value = 6.148826
$ws.Range("O6").Value = 0.9069818039462568
$ws.Range("P6").Value = 0.9069818039462569
$ws.Range("Q6").Value = 0.21164259092
$ws.Range("R6").Value = 1.90478331828
$ws.Range("S6").Value = 0.7679656018391702
$ws.Range("T6").Value = 0.7679656018391702

# --- Row 7 (new row) ---
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Nlgn1"
$ws.Range("C7").Value = "Nrxn2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.10326
$ws.Range("H7").Value = 0.30978
$ws.Range("I7").Value = 0.846726580949882
$ws.Range("J7").Value = 0.8467265809498818
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.206493
$ws.Range("N7").Value = 0.619479
$ws.Range("O7").Value = 0.09137617179715662
$ws.Range("P7").Value = 0.09137617179715661
$ws.Range("Q7").Value = 0.02132246718
$ws.Range("R7").Value = 0.19190220462
$ws.Range("S7").Value = 0.07737063352609545
$ws.Range("T7").Value = 0.07737063352609544
